$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment two id-like counters by 1
$ws.Range("B19").Value2 = 91806
$ws.Range("B21").Value2 = 79245

# Rows 28 and 29 hold two different bird observations that were recorded in
# the wrong row order. Swap their content back so each field lines up with
# the correct observation. Only the columns whose value actually differs
# between the two rows are touched (columns that already hold the same
# value in both rows - e.g. the shared location / observer / date fields -
# are left completely untouched).
$cols = @("A","B","E","F","G","L","M","Q","R","Z","AB")

foreach ($col in $cols) {
    $cellA = $ws.Range($col + "28")
    $cellB = $ws.Range($col + "29")
    $tmp = $cellA.Value2
    $cellA.Value2 = $cellB.Value2
    $cellB.Value2 = $tmp
}
